$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2903.6667
$ws.Range("I15").Value = 2903.6667
$ws.Range("K15").Value = 8711.000100000001
$ws.Range("M15").Value = -8542.000100000001

$ws.Range("H40").Value = 1911.4445
$ws.Range("J40").Value = 1880
$ws.Range("L40").Value = 1880
$ws.Range("N40").Value = -2230

$ws.Range("H80").Value = 675.3333
$ws.Range("I80").Value = 417.5
$ws.Range("J80").Value = 804.25
$ws.Range("K80").Value = 1252.5
$ws.Range("L80").Value = 2412.75
$ws.Range("M80").Value = -254.5
$ws.Range("N80").Value = -4408.75

$ws.Range("H83").Value = 675.3333
$ws.Range("I83").Value = 417.5
$ws.Range("J83").Value = 804.25
$ws.Range("K83").Value = 3757.5
$ws.Range("L83").Value = 7238.25
$ws.Range("M83").Value = 1234.5
$ws.Range("N83").Value = -17222.25

$ws.Range("H87").Value = 42872.5
$ws.Range("J87").Value = 42872.5
$ws.Range("L87").Value = 42872.5
$ws.Range("N87").Value = -45368.5

$ws.Range("H88").Value = 1546459.6
$ws.Range("I88").Value = 1499
$ws.Range("J88").Value = 2061446.5
$ws.Range("K88").Value = 1499
$ws.Range("L88").Value = 2061446.5
$ws.Range("M88").Value = -1093
$ws.Range("N88").Value = -2062258.5

$ws.Range("H90").Value = 42872.5
$ws.Range("J90").Value = 42872.5
$ws.Range("L90").Value = 128617.5
$ws.Range("N90").Value = -141097.5

$ws.Range("H91").Value = 1546459.6
$ws.Range("I91").Value = 1499
$ws.Range("J91").Value = 2061446.5
$ws.Range("K91").Value = 1499
$ws.Range("L91").Value = 2061446.5
$ws.Range("M91").Value = -95
$ws.Range("N91").Value = -2064254.5

$ws.Range("H113").Value = 2160
$ws.Range("I113").Value = 2007
$ws.Range("J113").Value = 2313
$ws.Range("K113").Value = 2007
$ws.Range("L113").Value = 2313
$ws.Range("M113").Value = 1247
$ws.Range("N113").Value = -8821

$ws.Range("H132").Value = 12830727
$ws.Range("I132").Value = 15159787
$ws.Range("K132").Value = 45479361
$ws.Range("M132").Value = -45476831

$ws.Range("H133").Value = 33770.4
$ws.Range("J133").Value = 33770.4
$ws.Range("L133").Value = 33770.4
$ws.Range("N133").Value = -43890.4

$ws.Range("H135").Value = 76925550
$ws.Range("I135").Value = 893.625
$ws.Range("J135").Value = 200005000
$ws.Range("K135").Value = 8042.625
$ws.Range("L135").Value = 1800045000
$ws.Range("M135").Value = -5507.625
$ws.Range("N135").Value = -1800050070

$ws.Range("H137").Value = 1406.8649
$ws.Range("I137").Value = 922.9524
$ws.Range("K137").Value = 2768.8572
$ws.Range("M137").Value = -218.8571999999999

$ws.Range("H138").Value = 440909.8
$ws.Range("J138").Value = 538661.9399999999
$ws.Range("L138").Value = 1615985.82
$ws.Range("N138").Value = -1626265.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4672.6
$ws.Range("I32").Value = 4582.3223
$ws.Range("K32").Value = 4582.3223
$ws.Range("M32").Value = -4295.3223

$ws.Range("H61").Value = 41667740
$ws.Range("I61").Value = 50000920
$ws.Range("K61").Value = 50000920
$ws.Range("M61").Value = -50000708

$ws.Range("H74").Value = 1314.2858
$ws.Range("I74").Value = 1314.2858
$ws.Range("K74").Value = 1314.2858
$ws.Range("M74").Value = -440.2858000000001

$ws.Range("H77").Value = 1314.2858
$ws.Range("I77").Value = 1314.2858
$ws.Range("K77").Value = 6571.429
$ws.Range("M77").Value = -2203.429

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 1598.1111
$ws.Range("I122").Value = 1415.6875
$ws.Range("J122").Value = 1863.4546
$ws.Range("K122").Value = 4247.0625
$ws.Range("L122").Value = 5590.3638
$ws.Range("M122").Value = -1797.0625
$ws.Range("N122").Value = -10490.3638

$ws.Range("H124").Value = 17495
$ws.Range("J124").Value = 17495
$ws.Range("L124").Value = 17495
$ws.Range("N124").Value = -27315

$ws.Range("H132").Value = 3054.348
$ws.Range("I132").Value = 2607.6
$ws.Range("K132").Value = 7822.799999999999
$ws.Range("M132").Value = -5292.799999999999

$ws.Range("H136").Value = 41667740
$ws.Range("I136").Value = 50000920
$ws.Range("K136").Value = 150002760
$ws.Range("M136").Value = -150000210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 50374.57
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 70414.39999999999
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 70414.39999999999
$ws.Range("M22").Value = 75
$ws.Range("N22").Value = -71114.39999999999

$ws.Range("H58").Value = 1618.6923
$ws.Range("I58").Value = 1368
$ws.Range("J58").Value = 2092.2222
$ws.Range("K58").Value = 1368
$ws.Range("L58").Value = 2092.2222
$ws.Range("M58").Value = -1165
$ws.Range("N58").Value = -2498.2222

$ws.Range("H99").Value = 1431.6154
$ws.Range("I99").Value = 1306.5555
$ws.Range("J99").Value = 1713
$ws.Range("K99").Value = 1306.5555
$ws.Range("L99").Value = 1713
$ws.Range("M99").Value = 191.4445000000001
$ws.Range("N99").Value = -4709

$ws.Range("H126").Value = 1431.6154
$ws.Range("I126").Value = 1306.5555
$ws.Range("J126").Value = 1713
$ws.Range("K126").Value = 3919.6665
$ws.Range("L126").Value = 5139
$ws.Range("M126").Value = -1449.6665
$ws.Range("N126").Value = -10079

$ws.Range("H132").Value = 2257.0833
$ws.Range("I132").Value = 2097.25
$ws.Range("J132").Value = 2576.75
$ws.Range("K132").Value = 6291.75
$ws.Range("L132").Value = 7730.25
$ws.Range("M132").Value = -3761.75
$ws.Range("N132").Value = -12790.25

$ws.Range("H134").Value = 20835684
$ws.Range("I134").Value = 2480.7368
$ws.Range("J134").Value = 100001860
$ws.Range("K134").Value = 7442.2104
$ws.Range("L134").Value = 300005580
$ws.Range("M134").Value = -4907.2104
$ws.Range("N134").Value = -300010650

$ws.Range("H136").Value = 1618.6923
$ws.Range("I136").Value = 1368
$ws.Range("J136").Value = 2092.2222
$ws.Range("K136").Value = 4104
$ws.Range("L136").Value = 6276.6666
$ws.Range("M136").Value = -1554
$ws.Range("N136").Value = -11376.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 408.69565
$ws.Range("I5").Value = 354.73685
$ws.Range("J5").Value = 665
$ws.Range("K5").Value = 1064.21055
$ws.Range("L5").Value = 1995
$ws.Range("M5").Value = -952.21055
$ws.Range("N5").Value = -2219

$ws.Range("H39").Value = 3400.125
$ws.Range("J39").Value = 3343
$ws.Range("L39").Value = 10029
$ws.Range("N39").Value = -10617

$ws.Range("H55").Value = 2158.6667
$ws.Range("J55").Value = 2499.6
$ws.Range("L55").Value = 7498.799999999999
$ws.Range("N55").Value = -7852.799999999999

$ws.Range("H113").Value = 670.1111
$ws.Range("I113").Value = 587.75
$ws.Range("J113").Value = 704.7895
$ws.Range("K113").Value = 1763.25
$ws.Range("L113").Value = 2114.3685
$ws.Range("M113").Value = 406.75
$ws.Range("N113").Value = -6454.3685

$ws.Range("H122").Value = 1160.3889
$ws.Range("J122").Value = 1189.8125
$ws.Range("L122").Value = 10708.3125
$ws.Range("N122").Value = -15608.3125

$ws.Range("H131").Value = 30304754
$ws.Range("J131").Value = 1867.7667
$ws.Range("L131").Value = 5603.300099999999
$ws.Range("N131").Value = -15683.3001

$ws.Range("H135").Value = 408.69565
$ws.Range("I135").Value = 354.73685
$ws.Range("J135").Value = 665
$ws.Range("K135").Value = 3192.63165
$ws.Range("L135").Value = 5985
$ws.Range("M135").Value = -657.6316500000003
$ws.Range("N135").Value = -11055

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1827.1818
$ws.Range("I7").Value = 1761.75
$ws.Range("J7").Value = 2001.6666
$ws.Range("K7").Value = 1761.75
$ws.Range("L7").Value = 2001.6666
$ws.Range("M7").Value = -1649.75
$ws.Range("N7").Value = -2225.6666

$ws.Range("H22").Value = 1006.75
$ws.Range("I22").Value = 1025.8
$ws.Range("J22").Value = 993.1429000000001
$ws.Range("K22").Value = 1025.8
$ws.Range("L22").Value = 993.1429000000001
$ws.Range("M22").Value = -730.8
$ws.Range("N22").Value = -1583.1429

$ws.Range("H27").Value = 1006.75
$ws.Range("I27").Value = 1025.8
$ws.Range("J27").Value = 993.1429000000001
$ws.Range("K27").Value = 1025.8
$ws.Range("L27").Value = 993.1429000000001
$ws.Range("M27").Value = -918.8
$ws.Range("N27").Value = -1207.1429

$ws.Range("H55").Value = 273
$ws.Range("J55").Value = 983.3333
$ws.Range("L55").Value = 983.3333
$ws.Range("N55").Value = -1329.3333

$ws.Range("H68").Value = 2098.75
$ws.Range("J68").Value = 1800
$ws.Range("L68").Value = 1800
$ws.Range("N68").Value = -3298

$ws.Range("H71").Value = 2098.75
$ws.Range("J71").Value = 1800
$ws.Range("L71").Value = 9000
$ws.Range("N71").Value = -16488

$ws.Range("H126").Value = 1827.1818
$ws.Range("I126").Value = 1761.75
$ws.Range("J126").Value = 2001.6666
$ws.Range("K126").Value = 5285.25
$ws.Range("L126").Value = 6004.9998
$ws.Range("M126").Value = -2815.25
$ws.Range("N126").Value = -10944.9998

$ws.Range("H132").Value = 2422.423
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 3817.1667
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 11451.5001
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -16511.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 9599.75
$ws.Range("J101").Value = 9599.75
$ws.Range("L101").Value = 9599.75
$ws.Range("N101").Value = -16089.75

$ws.Range("H122").Value = 16668249
$ws.Range("I122").Value = 27779786
$ws.Range("J122").Value = 941.6667
$ws.Range("K122").Value = 83339358
$ws.Range("L122").Value = 2825.0001
$ws.Range("M122").Value = -83336908
$ws.Range("N122").Value = -7725.0001

$ws.Range("H132").Value = 1927.762
$ws.Range("I132").Value = 1230.4375
$ws.Range("K132").Value = 3691.3125
$ws.Range("M132").Value = -1161.3125
